$wb = $excel.ActiveWorkbook

# --- Climates_to_remove: drop the "ET" climate row ---------------------
$wsClimates = $wb.Worksheets.Item("Climates_to_remove")
$wsClimates.Range("A2").EntireRow.Delete()
$wsClimates.Range("E9").Select()

# --- Other settings: tweak the region / host-table flags ----------------
$wsOther = $wb.Worksheets.Item("Other settings")
$wsOther.Range("B2").Value = "yes"
$wsOther.Range("B3").Value = "Global"
$wsOther.Range("C12").Select()

# --- tech: adjust Europe bounding box, add Peru as a new region ---------
$wsTech = $wb.Worksheets.Item("tech")
$wsTech.Range("B7").Value = -30
$wsTech.Range("D7").Value = 35

$wsTech.Range("A10").Value = "Peru"
$wsTech.Range("B10").Value = -85
$wsTech.Range("C10").Value = -65
$wsTech.Range("D10").Value = -20
$wsTech.Range("E10").Value = 5
$wsTech.Range("F10").Value = 5
$wsTech.Range("G10").Value = 5
$wsTech.Range("H39").Select()

# --- Pest_list: start tracking the new pest ------------------------------
$wsPest = $wb.Worksheets.Item("Pest_list")
$wsPest.Range("A2").Value = "Elasmopalpus lignosellus"
$wsPest.Activate()
$wsPest.Range("A3").Select()
